# Scheduled market-data refresh: update currentAveragePrice(NQ/HQ) and
# derived Leve profit columns (H:N) across the per-job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 9449.666999999999
$ws.Cells.Item(43, 9).Value = 4949.5
$ws.Cells.Item(43, 10).Value = 10735.429
$ws.Cells.Item(43, 11).Value = 4949.5
$ws.Cells.Item(43, 12).Value = 10735.429
$ws.Cells.Item(43, 13).Value = -4880.5
$ws.Cells.Item(43, 14).Value = -10873.429
$ws.Cells.Item(62, 8).Value = 7018.5557
$ws.Cells.Item(62, 10).Value = 7065.4287
$ws.Cells.Item(62, 12).Value = 7065.4287
$ws.Cells.Item(62, 14).Value = -8313.4287
$ws.Cells.Item(64, 8).Value = 7428.2856
$ws.Cells.Item(64, 10).Value = 7999.6665
$ws.Cells.Item(64, 12).Value = 7999.6665
$ws.Cells.Item(64, 14).Value = -8495.666499999999
$ws.Cells.Item(65, 8).Value = 7018.5557
$ws.Cells.Item(65, 10).Value = 7065.4287
$ws.Cells.Item(65, 12).Value = 35327.14350000001
$ws.Cells.Item(65, 14).Value = -41567.14350000001
$ws.Cells.Item(67, 8).Value = 7428.2856
$ws.Cells.Item(67, 10).Value = 7999.6665
$ws.Cells.Item(67, 12).Value = 7999.6665
$ws.Cells.Item(67, 14).Value = -9715.666499999999
$ws.Cells.Item(137, 8).Value = 1617.1464
$ws.Cells.Item(137, 9).Value = 1676.0769
$ws.Cells.Item(137, 11).Value = 5028.2307
$ws.Cells.Item(137, 13).Value = -2478.2307
$ws.Cells.Item(138, 8).Value = 269961.28
$ws.Cells.Item(138, 9).Value = 78149.766
$ws.Cells.Item(138, 10).Value = 359016.66
$ws.Cells.Item(138, 11).Value = 234449.298
$ws.Cells.Item(138, 12).Value = 1077049.98
$ws.Cells.Item(138, 13).Value = -229309.298
$ws.Cells.Item(138, 14).Value = -1087329.98

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1066.2646
$ws.Cells.Item(2, 9).Value = 1053.7693
$ws.Cells.Item(2, 11).Value = 1053.7693
$ws.Cells.Item(2, 13).Value = -940.7692999999999
$ws.Cells.Item(32, 8).Value = 8647.309999999999
$ws.Cells.Item(32, 9).Value = 9153.865
$ws.Cells.Item(32, 11).Value = 9153.865
$ws.Cells.Item(32, 13).Value = -8866.865
$ws.Cells.Item(43, 8).Value = 14418.286
$ws.Cells.Item(43, 9).Value = 11088
$ws.Cells.Item(43, 10).Value = 15750.4
$ws.Cells.Item(43, 11).Value = 11088
$ws.Cells.Item(43, 12).Value = 15750.4
$ws.Cells.Item(43, 13).Value = -10775
$ws.Cells.Item(43, 14).Value = -16376.4
$ws.Cells.Item(61, 8).Value = 4488.8086
$ws.Cells.Item(61, 9).Value = 2880.186
$ws.Cells.Item(61, 11).Value = 2880.186
$ws.Cells.Item(61, 13).Value = -2668.186
$ws.Cells.Item(74, 8).Value = 3010.375
$ws.Cells.Item(74, 9).Value = 2451.818
$ws.Cells.Item(74, 11).Value = 2451.818
$ws.Cells.Item(74, 13).Value = -1577.818
$ws.Cells.Item(77, 8).Value = 3010.375
$ws.Cells.Item(77, 9).Value = 2451.818
$ws.Cells.Item(77, 11).Value = 12259.09
$ws.Cells.Item(77, 13).Value = -7891.09
$ws.Cells.Item(116, 8).Value = 1066.2646
$ws.Cells.Item(116, 9).Value = 1053.7693
$ws.Cells.Item(116, 11).Value = 1053.7693
$ws.Cells.Item(116, 13).Value = 1240.2307
$ws.Cells.Item(132, 8).Value = 2346.95
$ws.Cells.Item(132, 9).Value = 2196.849
$ws.Cells.Item(132, 11).Value = 6590.547
$ws.Cells.Item(132, 13).Value = -4060.547
$ws.Cells.Item(136, 8).Value = 4488.8086
$ws.Cells.Item(136, 9).Value = 2880.186
$ws.Cells.Item(136, 11).Value = 8640.558000000001
$ws.Cells.Item(136, 13).Value = -6090.558000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1066.2646
$ws.Cells.Item(3, 9).Value = 1053.7693
$ws.Cells.Item(3, 11).Value = 1053.7693
$ws.Cells.Item(3, 13).Value = -939.7692999999999
$ws.Cells.Item(100, 8).Value = 4950
$ws.Cells.Item(100, 10).Value = 4950
$ws.Cells.Item(100, 12).Value = 4950
$ws.Cells.Item(100, 14).Value = -7114
$ws.Cells.Item(107, 8).Value = 2077.2856
$ws.Cells.Item(107, 9).Value = 2131.6667
$ws.Cells.Item(107, 10).Value = 2036.5
$ws.Cells.Item(107, 11).Value = 2131.6667
$ws.Cells.Item(107, 12).Value = 2036.5
$ws.Cells.Item(107, 13).Value = -211.6667000000002
$ws.Cells.Item(107, 14).Value = -5876.5
$ws.Cells.Item(134, 8).Value = 9279.4375
$ws.Cells.Item(134, 9).Value = 2314
$ws.Cells.Item(134, 11).Value = 6942
$ws.Cells.Item(134, 13).Value = -4407

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1938.1311
$ws.Cells.Item(31, 10).Value = 3584.55
$ws.Cells.Item(31, 12).Value = 3584.55
$ws.Cells.Item(31, 14).Value = -4174.55
$ws.Cells.Item(34, 8).Value = 1938.1311
$ws.Cells.Item(34, 10).Value = 3584.55
$ws.Cells.Item(34, 12).Value = 3584.55
$ws.Cells.Item(34, 14).Value = -3988.55
$ws.Cells.Item(122, 8).Value = 2767.2646
$ws.Cells.Item(122, 9).Value = 2608.074
$ws.Cells.Item(122, 10).Value = 3381.2856
$ws.Cells.Item(122, 11).Value = 7824.222
$ws.Cells.Item(122, 12).Value = 10143.8568
$ws.Cells.Item(122, 13).Value = -5374.222
$ws.Cells.Item(122, 14).Value = -15043.8568
$ws.Cells.Item(132, 8).Value = 1430969.4
$ws.Cells.Item(132, 9).Value = 1602386.6
$ws.Cells.Item(132, 11).Value = 4807159.800000001
$ws.Cells.Item(132, 13).Value = -4804629.800000001
$ws.Cells.Item(140, 8).Value = 93077.5
$ws.Cells.Item(140, 10).Value = 107539.57
$ws.Cells.Item(140, 12).Value = 107539.57
$ws.Cells.Item(140, 14).Value = -117899.57
$ws.Cells.Item(141, 8).Value = 239340.11
$ws.Cells.Item(141, 10).Value = 239340.11
$ws.Cells.Item(141, 12).Value = 239340.11
$ws.Cells.Item(141, 14).Value = -249700.11

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(99, 8).Value = 39138.75
$ws.Cells.Item(99, 10).Value = 40444.285
$ws.Cells.Item(99, 12).Value = 40444.285
$ws.Cells.Item(99, 14).Value = -44936.285
$ws.Cells.Item(102, 8).Value = 28143.29
$ws.Cells.Item(102, 9).Value = 1469.72
$ws.Cells.Item(102, 11).Value = 1469.72
$ws.Cells.Item(102, 13).Value = 152.28
$ws.Cells.Item(113, 8).Value = 1897.3
$ws.Cells.Item(113, 9).Value = 1963.6666
$ws.Cells.Item(113, 10).Value = 1300
$ws.Cells.Item(113, 11).Value = 1963.6666
$ws.Cells.Item(113, 12).Value = 1300
$ws.Cells.Item(113, 13).Value = 206.3334
$ws.Cells.Item(132, 8).Value = 9807240
$ws.Cells.Item(132, 9).Value = 13336273
$ws.Cells.Item(132, 11).Value = 40008819
$ws.Cells.Item(132, 13).Value = -40006289
$ws.Cells.Item(113, 14).Value = -5640

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 3919.75
$ws.Cells.Item(16, 10).Value = 4395
$ws.Cells.Item(16, 12).Value = 4395
$ws.Cells.Item(16, 14).Value = -4735
$ws.Cells.Item(95, 8).Value = 48926
$ws.Cells.Item(95, 10).Value = 48926
$ws.Cells.Item(95, 12).Value = 48926
$ws.Cells.Item(95, 14).Value = -54418
$ws.Cells.Item(97, 8).Value = 63937.4
$ws.Cells.Item(97, 10).Value = 63937.4
$ws.Cells.Item(97, 12).Value = 63937.4
$ws.Cells.Item(97, 14).Value = -65919.39999999999
$ws.Cells.Item(132, 8).Value = 3118.5557
$ws.Cells.Item(132, 9).Value = 3105.861
$ws.Cells.Item(132, 11).Value = 9317.582999999999
$ws.Cells.Item(132, 13).Value = -6787.582999999999
$ws.Cells.Item(134, 8).Value = 83987
$ws.Cells.Item(134, 10).Value = 87076.73
$ws.Cells.Item(134, 12).Value = 87076.73
$ws.Cells.Item(134, 14).Value = -97216.73

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(117, 8).Value = 89897
$ws.Cells.Item(117, 10).Value = 89897
$ws.Cells.Item(117, 12).Value = 89897
$ws.Cells.Item(117, 14).Value = -99075
$ws.Cells.Item(136, 8).Value = 10099.956
$ws.Cells.Item(136, 9).Value = 24255.445
$ws.Cells.Item(136, 11).Value = 72766.33499999999
$ws.Cells.Item(136, 13).Value = -70216.33499999999
